$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row strings: "<field>_old" -> "<field>_FV2310" and
#    "<field>_new" -> "<field>_FV2404" (columns A:J and L:U; column K stays "diff").
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the data range into an Excel Table (ListObject) so the header row
#    gets an AutoFilter + banded rows, matching xl/tables/table1.xml.
$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.TableStyle = ""

# 3. Freeze the header row (pane split after row 1).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "edit applied"
